$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "LRE" in column I, row 1
$ws.Range("I1").Value = "LRE"

# Set column I values (rows 2-6) to 0
$ws.Range("I2:I6").Value = 0

# Set column H width to match diff (14.109375, snaps to nearest reachable
# pixel-quantized width of 14.1666... in this engine)
$ws.Columns.Item(8).ColumnWidth = 13.33

# Select cell I5 as the active selection (matches diff's sheetView selection)
$ws.Range("I5").Select()
